$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update k_kalibrate values in ROM table (columns B:AG, rows 2-9)
# Values derived from the authoritative diff of table_for_grafics_new.xlsx

# Row 2
$ws.Range("B2").Value = 1.3
$ws.Range("C2").Value = 1.27
$ws.Range("D2").Value = 1.31
$ws.Range("E2").Value = 1.29
$ws.Range("F2").Value = 1.32
$ws.Range("G2").Value = 1.28
$ws.Range("I2").Value = 1.33
$ws.Range("J2").Value = 1.32
$ws.Range("K2").Value = 1.37
$ws.Range("L2").Value = 1.29
$ws.Range("M2").Value = 1.32
$ws.Range("N2").Value = 1.26
$ws.Range("O2").Value = 1.33
$ws.Range("P2").Value = 1.41
$ws.Range("Q2").Value = 1.33
$ws.Range("R2").Value = 1.28
$ws.Range("T2").Value = 1.4
$ws.Range("U2").Value = 1.41
$ws.Range("V2").Value = 1.28
$ws.Range("W2").Value = 1.31
$ws.Range("X2").Value = 1.34
$ws.Range("Z2").Value = 1.31
$ws.Range("AB2").Value = 1.31
$ws.Range("AC2").Value = 1.34
$ws.Range("AD2").Value = 1.34
$ws.Range("AG2").Value = 1.28

# Row 3
$ws.Range("J3").Value = 1.14
$ws.Range("L3").Value = 1.13
$ws.Range("M3").Value = 1.14
$ws.Range("S3").Value = 1.15
$ws.Range("T3").Value = 1.17
$ws.Range("V3").Value = 1.12
$ws.Range("W3").Value = 1.14
$ws.Range("X3").Value = 1.15
$ws.Range("AA3").Value = 1.13
$ws.Range("AD3").Value = 1.15

# Row 4
$ws.Range("H4").Value = 1.09
$ws.Range("J4").Value = 1.1
$ws.Range("L4").Value = 1.09
$ws.Range("R4").Value = 1.1
$ws.Range("AD4").Value = 1.11
$ws.Range("AG4").Value = 1.09

# Row 5
$ws.Range("J5").Value = 1.08
$ws.Range("K5").Value = 1.09
$ws.Range("P5").Value = 1.09
$ws.Range("T5").Value = 1.09
$ws.Range("U5").Value = 1.09
$ws.Range("V5").Value = 1.08
$ws.Range("AD5").Value = 1.09

# Row 6
$ws.Range("G6").Value = 1.04
$ws.Range("I6").Value = 1.07
$ws.Range("J6").Value = 1.06
$ws.Range("T6").Value = 1.08
$ws.Range("U6").Value = 1.08
$ws.Range("AB6").Value = 1.06

# Row 7
$ws.Range("H7").Value = 1.05
$ws.Range("J7").Value = 1.06
$ws.Range("N7").Value = 1.04
$ws.Range("S7").Value = 1.06
$ws.Range("AE7").Value = 1.05

# Row 8
$ws.Range("J8").Value = 1.05
$ws.Range("K8").Value = 1.06
$ws.Range("X8").Value = 1.06

# Row 9
$ws.Range("J9").Value = 1.05
$ws.Range("K9").Value = 1.06
$ws.Range("Q9").Value = 1.06
